$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (A2) used to hold "fsssacct3@gmail.com"; the account used for
# Google sign-in testing is now the one that used to be on row 3
# ("fsqa1tpn@gmail.com"). Overwrite A2's text in place - this keeps the
# existing hyperlink attached to A2 untouched (same r:id/target as before),
# exactly like a user typing over the cell text in Excel.
$ws.Range("A2").Value = "fsqa1tpn@gmail.com"

# Row 3 is no longer needed - clear its two cells (keeps the row/cell
# styling, just drops the content), matching the emptied <c .../> cells.
$ws.Range("A3:B3").ClearContents()

# Remove the (now stale) hyperlinks that used to target A3/B3. This engine
# only supports clearing the *entire* hyperlink collection in one shot, so
# wipe it and re-create just the two links that must survive (A2 -> the
# original rId1 mail target, B2 -> the original rId2 mail target).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:fsssacct3@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:Admin@123")

# Re-adding hyperlinks stamps a brand new cell style; put A2:B2 back on the
# workbook's built-in "Hyperlink" style so they render/save like the
# original (underlined, themed) cells instead of a duplicate style entry.
$ws.Range("A2:B2").Style = "Hyperlink"

# Match the saved selection (A2 instead of the old D3).
$ws.Range("A2").Select() | Out-Null
